$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($CellRange, $Text)
    $CellRange.Value = "'" + $Text
    $CellRange.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "26.401.99"
Set-TextCell $ws.Range("E2") "  +1.54%  "
Set-TextCell $ws.Range("D3") "1.667.15"
Set-TextCell $ws.Range("E3") "  +1.18%  "
Set-TextCell $ws.Range("D4") "1.004"
Set-TextCell $ws.Range("E4") "  +0.32%  "
Set-TextCell $ws.Range("D5") "219.71"
Set-TextCell $ws.Range("E5") "  +2.61%  "
Set-TextCell $ws.Range("D6") "0.5248"
Set-TextCell $ws.Range("E6") "  +0.43%  "
Set-TextCell $ws.Range("E7") "  +0.30%  "
Set-TextCell $ws.Range("D8") "0.2668"
Set-TextCell $ws.Range("E8") "  +1.93%  "
Set-TextCell $ws.Range("D9") "0.06351"
Set-TextCell $ws.Range("E9") "  -0.03%  "
Set-TextCell $ws.Range("D10") "21.63"
Set-TextCell $ws.Range("E10") "  +4.25%  "
Set-TextCell $ws.Range("D11") "0.07767"
Set-TextCell $ws.Range("E11") "  +0.83%  "
Set-TextCell $ws.Range("B12") "WrappedEther"
Set-TextCell $ws.Range("C12") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws.Range("D12") "1.679.47"
Set-TextCell $ws.Range("E12") "  +1.84%  "
Set-TextCell $ws.Range("B13") "Polkadot"
Set-TextCell $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws.Range("D13") "4.459"
Set-TextCell $ws.Range("E13") "  +0.57%  "
Set-TextCell $ws.Range("D14") "0.5512"
Set-TextCell $ws.Range("E14") "  -0.03%  "
Set-TextCell $ws.Range("D15") "0.0₅8252"
Set-TextCell $ws.Range("E15") "  +0.20%  "
Set-TextCell $ws.Range("D16") "65.43"
Set-TextCell $ws.Range("E16") "  +1.10%  "
Set-TextCell $ws.Range("D17") "26.425.88"
Set-TextCell $ws.Range("E17") "  +1.51%  "
Set-TextCell $ws.Range("D18") "1.003"
Set-TextCell $ws.Range("E18") "  +0.25%  "
Set-TextCell $ws.Range("D19") "4.727"
Set-TextCell $ws.Range("E19") "  -0.09%  "
Set-TextCell $ws.Range("D20") "193.29"
Set-TextCell $ws.Range("E20") "  +1.60%  "
Set-TextCell $ws.Range("D21") "10.25"
Set-TextCell $ws.Range("E21") "  +0.35%  "
Set-TextCell $ws.Range("D22") "6.248"
Set-TextCell $ws.Range("E22") "  -1.07%  "
Set-TextCell $ws.Range("D23") "1.006"
Set-TextCell $ws.Range("E23") "  +0.29%  "
Set-TextCell $ws.Range("B24") "Monero"
Set-TextCell $ws.Range("C24") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D24") "139.24"
Set-TextCell $ws.Range("E24") "  -2.69%  "
Set-TextCell $ws.Range("B25") "Stellar"
Set-TextCell $ws.Range("C25") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D25") "0.1262"
Set-TextCell $ws.Range("E25") "  +1.32%  "
Set-TextCell $ws.Range("D26") "7.362"
Set-TextCell $ws.Range("E26") "  -0.54%  "
Set-TextCell $ws.Range("D27") "16.20"
Set-TextCell $ws.Range("E27") "  +1.24%  "
Set-TextCell $ws.Range("D28") "1.413"
Set-TextCell $ws.Range("E28") "  +0.37%  "
Set-TextCell $ws.Range("D29") "0.06126"
Set-TextCell $ws.Range("E29") "  +3.10%  "
Set-TextCell $ws.Range("D30") "1.289"
Set-TextCell $ws.Range("E30") "  +2.40%  "
Set-TextCell $ws.Range("D31") "3.584"
Set-TextCell $ws.Range("E31") "  +4.59%  "
Set-TextCell $ws.Range("D32") "3.388"
Set-TextCell $ws.Range("E32") "  -0.50%  "
Set-TextCell $ws.Range("D33") "1.673"
Set-TextCell $ws.Range("E33") "  +1.82%  "
Set-TextCell $ws.Range("D34") "1.000"
Set-TextCell $ws.Range("E34") "  +0.71%  "
Set-TextCell $ws.Range("D35") "2.420"
Set-TextCell $ws.Range("E35") "  +1.02%  "
Set-TextCell $ws.Range("D36") "0.6031"
Set-TextCell $ws.Range("E36") "  +6.95%  "
Set-TextCell $ws.Range("D37") "2.775"
Set-TextCell $ws.Range("E37") "  +0.65%  "
Set-TextCell $ws.Range("D38") "0.01609"
Set-TextCell $ws.Range("E38") "  +0.08%  "
Set-TextCell $ws.Range("D39") "6.019"
Set-TextCell $ws.Range("E39") "  +2.73%  "
Set-TextCell $ws.Range("D40") "1.085.93"
Set-TextCell $ws.Range("E40") "  +5.91%  "
Set-TextCell $ws.Range("D41") "0.8557"
Set-TextCell $ws.Range("E41") "  -0.18%  "
Set-TextCell $ws.Range("D42") "1.002"
Set-TextCell $ws.Range("E42") "  +0.21%  "
Set-TextCell $ws.Range("D43") "100.52"
Set-TextCell $ws.Range("E43") "  +1.55%  "
Set-TextCell $ws.Range("D44") "1.811.95"
Set-TextCell $ws.Range("E44") "  +0.83%  "
Set-TextCell $ws.Range("B45") "BabyDogeCoin"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws.Range("D45") "0.0₈109"
Set-TextCell $ws.Range("E45") "  +1.16%  "
Set-TextCell $ws.Range("B46") "Aave"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Range("D46") "57.79"
Set-TextCell $ws.Range("E46") "  +3.76%  "
Set-TextCell $ws.Range("B47") "EnergySwap"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D47") "8.165"
Set-TextCell $ws.Range("E47") "  +0.86%  "
Set-TextCell $ws.Range("B48") "Frax"
Set-TextCell $ws.Range("C48") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws.Range("D48") "1.004"
Set-TextCell $ws.Range("E48") "  +0.17%  "
Set-TextCell $ws.Range("D49") "0.05206"
Set-TextCell $ws.Range("E49") "  +1.17%  "
Set-TextCell $ws.Range("D50") "1.477"
Set-TextCell $ws.Range("E50") "  +6.80%  "
Set-TextCell $ws.Range("D51") "0.4231"
